$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5235.5
$ws.Range("I76").Value = 2765
$ws.Range("K76").Value = 2765
$ws.Range("M76").Value = -2450
$ws.Range("H79").Value = 5235.5
$ws.Range("I79").Value = 2765
$ws.Range("K79").Value = 2765
$ws.Range("M79").Value = -1673
$ws.Range("H137").Value = 3012.5908
$ws.Range("I137").Value = 1785.6316
$ws.Range("K137").Value = 5356.8948
$ws.Range("M137").Value = -2806.8948

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2847.4333
$ws.Range("I61").Value = 2165.8125
$ws.Range("K61").Value = 2165.8125
$ws.Range("M61").Value = -1953.8125
$ws.Range("H74").Value = 1755.5555
$ws.Range("I74").Value = 1486.5714
$ws.Range("J74").Value = 2697
$ws.Range("K74").Value = 1486.5714
$ws.Range("L74").Value = 2697
$ws.Range("M74").Value = -612.5714
$ws.Range("N74").Value = -4445
$ws.Range("H77").Value = 1755.5555
$ws.Range("I77").Value = 1486.5714
$ws.Range("J77").Value = 2697
$ws.Range("K77").Value = 7432.857
$ws.Range("L77").Value = 13485
$ws.Range("M77").Value = -3064.857
$ws.Range("N77").Value = -22221
$ws.Range("H132").Value = 32443.518
$ws.Range("I132").Value = 1723.875
$ws.Range("J132").Value = 179897.8
$ws.Range("K132").Value = 5171.625
$ws.Range("L132").Value = 539693.3999999999
$ws.Range("M132").Value = -2641.625
$ws.Range("N132").Value = -544753.3999999999
$ws.Range("H136").Value = 2847.4333
$ws.Range("I136").Value = 2165.8125
$ws.Range("K136").Value = 6497.4375
$ws.Range("M136").Value = -3947.4375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2609.658
$ws.Range("I86").Value = 1862.9642
$ws.Range("J86").Value = 4700.4
$ws.Range("K86").Value = 1862.9642
$ws.Range("L86").Value = 4700.4
$ws.Range("M86").Value = -739.9641999999999
$ws.Range("N86").Value = -6946.4
$ws.Range("H89").Value = 2609.658
$ws.Range("I89").Value = 1862.9642
$ws.Range("J89").Value = 4700.4
$ws.Range("K89").Value = 9314.821
$ws.Range("L89").Value = 23502
$ws.Range("M89").Value = -3698.821
$ws.Range("N89").Value = -34734
$ws.Range("H134").Value = 10322.72
$ws.Range("I134").Value = 3908.1667
$ws.Range("K134").Value = 11724.5001
$ws.Range("M134").Value = -9189.500100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7299
$ws.Range("I16").Value = 7643
$ws.Range("K16").Value = 7643
$ws.Range("M16").Value = -7356
$ws.Range("H105").Value = 4060.5
$ws.Range("I105").Value = 2299.5833
$ws.Range("J105").Value = 7582.3335
$ws.Range("K105").Value = 2299.5833
$ws.Range("L105").Value = 7582.3335
$ws.Range("M105").Value = -552.5832999999998
$ws.Range("N105").Value = -11076.3335
$ws.Range("H113").Value = 7299
$ws.Range("I113").Value = 7643
$ws.Range("K113").Value = 7643
$ws.Range("M113").Value = -5473
$ws.Range("H132").Value = 4469.1177
$ws.Range("I132").Value = 3266.074
$ws.Range("K132").Value = 9798.222
$ws.Range("M132").Value = -7268.222
$ws.Range("H141").Value = 59616.75
$ws.Range("I141").Value = 30000
$ws.Range("K141").Value = 30000
$ws.Range("M141").Value = -24820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 14596.96
$ws.Range("I3").Value = 3368.2666
$ws.Range("J3").Value = 31440
$ws.Range("K3").Value = 10104.7998
$ws.Range("L3").Value = 94320
$ws.Range("M3").Value = -9992.799800000001
$ws.Range("N3").Value = -94544
$ws.Range("H14").Value = 1643.125
$ws.Range("I14").Value = 1643.125
$ws.Range("K14").Value = 4929.375
$ws.Range("M14").Value = -4756.375
$ws.Range("H108").Value = 4648.5
$ws.Range("I108").Value = 578.2
$ws.Range("K108").Value = 1734.6
$ws.Range("M108").Value = 1145.4
$ws.Range("H112").Value = 14481.9
$ws.Range("I112").Value = 2139.6667
$ws.Range("K112").Value = 6419.000100000001
$ws.Range("M112").Value = -5311.000100000001
$ws.Range("H131").Value = 36112988
$ws.Range("J131").Value = 19610198
$ws.Range("L131").Value = 58830594
$ws.Range("N131").Value = -58840674
$ws.Range("H132").Value = 4482.5713
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 4979.6665
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 44816.9985
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -49876.9985
$ws.Range("H134").Value = 4907.069
$ws.Range("I134").Value = 4907.069
$ws.Range("K134").Value = 14721.207
$ws.Range("M134").Value = -9651.207000000002
$ws.Range("H140").Value = 49023604
$ws.Range("I140").Value = 64106252
$ws.Range("K140").Value = 192318756
$ws.Range("M140").Value = -192313576

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4885
$ws.Range("I132").Value = 4143.385
$ws.Range("K132").Value = 12430.155
$ws.Range("M132").Value = -9900.155000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4692.3228
$ws.Range("I7").Value = 3097
$ws.Range("K7").Value = 3097
$ws.Range("M7").Value = -2985
$ws.Range("H47").Value = 32499.166
$ws.Range("I47").Value = 26999.2
$ws.Range("J47").Value = 59999
$ws.Range("K47").Value = 26999.2
$ws.Range("L47").Value = 59999
$ws.Range("M47").Value = -26509.2
$ws.Range("N47").Value = -60979
$ws.Range("H52").Value = 32499.166
$ws.Range("I52").Value = 26999.2
$ws.Range("J52").Value = 59999
$ws.Range("K52").Value = 26999.2
$ws.Range("L52").Value = 59999
$ws.Range("M52").Value = -26766.2
$ws.Range("N52").Value = -60465
$ws.Range("H93").Value = 3031.9
$ws.Range("I93").Value = 3868.077
$ws.Range("J93").Value = 2392.4707
$ws.Range("K93").Value = 3868.077
$ws.Range("L93").Value = 2392.4707
$ws.Range("M93").Value = -2620.077
$ws.Range("N93").Value = -4888.4707
$ws.Range("H122").Value = 6986.1665
$ws.Range("I122").Value = 2776.5
$ws.Range("J122").Value = 7828.1
$ws.Range("K122").Value = 8329.5
$ws.Range("L122").Value = 23484.3
$ws.Range("M122").Value = -5879.5
$ws.Range("N122").Value = -28384.3
$ws.Range("H126").Value = 4692.3228
$ws.Range("I126").Value = 3097
$ws.Range("K126").Value = 9291
$ws.Range("M126").Value = -6821

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 200
$ws.Range("I14").Value = 200
$ws.Range("K14").Value = 200
$ws.Range("M14").Value = -32
$ws.Range("H81").Value = 789.3333
$ws.Range("I81").Value = 789.3333
$ws.Range("K81").Value = 1578.6666
$ws.Range("M81").Value = -517.6666
$ws.Range("H84").Value = 789.3333
$ws.Range("I84").Value = 789.3333
$ws.Range("K84").Value = 7893.333000000001
$ws.Range("M84").Value = -2589.333000000001
$ws.Range("H100").Value = 1386.2727
$ws.Range("I100").Value = 1143.1428
$ws.Range("J100").Value = 1811.75
$ws.Range("K100").Value = 2286.2856
$ws.Range("L100").Value = 3623.5
$ws.Range("M100").Value = -1745.2856
$ws.Range("N100").Value = -4705.5
$ws.Range("H139").Value = 81921.336
$ws.Range("J139").Value = 81921.336
$ws.Range("L139").Value = 81921.336
$ws.Range("N139").Value = -92201.336
